$d = $word.ActiveDocument

$replacements = @(
    @("2024-11-28 Thursday", "2024-11-29 Friday"),
    @("63×57=3591", "42×30=1260"),
    @("64×81=5184", "19×35=665"),
    @("21×57=1197", "35×58=2030"),
    @("60×59=3540", "39×73=2847"),
    @("86×96=8256", "58×37=2146"),
    @("53×37=1961", "90×38=3420"),
    @("97×67=6499", "38×80=3040"),
    @("67×67=4489", "95×56=5320"),
    @("71×60=4260", "74×22=1628"),
    @("11×28=308", "78×94=7332"),
    @("96×18=1728", "90×33=2970"),
    @("48×34=1632", "94×41=3854"),
    @("91×25=2275", "72×81=5832"),
    @("47×68=3196", "81×14=1134"),
    @("45×48=2160", "54×91=4914"),
    @("90×92=8280", "31×69=2139"),
    @("52×25=1300", "57×22=1254"),
    @("31×99=3069", "90×97=8730"),
    @("68×84=5712", "55×67=3685"),
    @("21×58=1218", "85×29=2465"),
    @("37×47=1739", "70×60=4200"),
    @("66×11=726", "52×76=3952"),
    @("72×42=3024", "87×25=2175"),
    @("21×20=420", "11×16=176"),
    @("37×90=3330", "19×82=1558")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
